$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bulk-updated
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2..491).
$ws.Range("C2:C491").Value = 45188
